$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellAddr, [string]$val) {
    $c = $ws.Range($cellAddr)
    # Leading apostrophe forces Excel to treat the value as literal text
    # (same as the quote-prefix trick in the Excel UI), preventing
    # numeric-looking strings (e.g. '591.56', '1.07') from being
    # silently converted to numbers.
    $c.Value = [string]::Concat("'", $val)
    # Resetting the style back to "Normal" clears the quotePrefix/text
    # number-format flag Excel attaches to the cell when it applies the
    # quote-prefix, so the cell keeps its original (default) style.
    $c.Style = "Normal"
}

Set-TextValue "D2" '63.597.07'
Set-TextValue "E2" '  -1.59%  '
Set-TextValue "D3" '3.074.16'
Set-TextValue "E3" '  -1.48%  '
Set-TextValue "E4" '  -0.79%  '
Set-TextValue "D5" '591.56'
Set-TextValue "E5" '  +0.32%  '
Set-TextValue "D6" '154.75'
Set-TextValue "E6" '  +1.45%  '
Set-TextValue "E7" '  -0.31%  '
Set-TextValue "D8" '0.535'
Set-TextValue "E8" '  +0.81%  '
Set-TextValue "D9" '3.072.98'
Set-TextValue "E9" '  -1.47%  '
Set-TextValue "E10" '  -1.57%  '
Set-TextValue "D11" '5.91'
Set-TextValue "E11" '  -0.44%  '
Set-TextValue "D12" '0.450'
Set-TextValue "E12" '  -2.50%  '
Set-TextValue "E13" '  -3.18%  '
Set-TextValue "D14" '36.52'
Set-TextValue "E14" '  -3.66%  '
Set-TextValue "E15" '  +0.55%  '
Set-TextValue "D16" '3.580.18'
Set-TextValue "E16" '  -1.75%  '
Set-TextValue "D17" '7.17'
Set-TextValue "E17" '  -1.08%  '
Set-TextValue "D18" '63.493.38'
Set-TextValue "E18" '  -1.09%  '
Set-TextValue "D19" '3.073.25'
Set-TextValue "E19" '  -1.62%  '
Set-TextValue "D20" '481.10'
Set-TextValue "E20" '  +1.93%  '
Set-TextValue "D21" '14.41'
Set-TextValue "E21" '  -3.59%  '
Set-TextValue "E22" '  -4.34%  '
Set-TextValue "D23" '7.54'
Set-TextValue "E23" '  -1.00%  '
Set-TextValue "E24" '  +0.69%  '
Set-TextValue "D25" '81.67'
Set-TextValue "E25" '  -0.25%  '
Set-TextValue "D26" '12.83'
Set-TextValue "E26" '  -3.63%  '
Set-TextValue "D27" '10.64'
Set-TextValue "E27" '  +7.82%  '
Set-TextValue "E28" '  +0.13%  '
Set-TextValue "D29" '7.60'
Set-TextValue "E29" '  +2.96%  '
Set-TextValue "E30" '  -0.88%  '
Set-TextValue "E31" '  -0.67%  '
Set-TextValue "D32" '2.19'
Set-TextValue "E32" '  -1.45%  '
Set-TextValue "E33" '  -4.72%  '
Set-TextValue "D34" '27.17'
Set-TextValue "E34" '  -1.32%  '
Set-TextValue "B35" 'Mantle'
Set-TextValue "C35" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D35" '1.07'
Set-TextValue "E35" '  +0.79%  '
Set-TextValue "B36" 'PEPE'
Set-TextValue "C36" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D36" [string]::Concat('0.0', [char]8323, '0829')
Set-TextValue "E36" '  -3.50%  '
Set-TextValue "B37" 'Filecoin'
Set-TextValue "C37" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D37" '6.07'
Set-TextValue "E37" '  -1.48%  '
Set-TextValue "B38" 'dogwifhat'
Set-TextValue "C38" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D38" '3.29'
Set-TextValue "E38" '  -1.87%  '
Set-TextValue "E39" '  -2.73%  '
Set-TextValue "D40" '50.58'
Set-TextValue "E40" '  -0.85%  '
Set-TextValue "D41" '9.19'
Set-TextValue "E41" '  -1.66%  '
Set-TextValue "D42" '441.11'
Set-TextValue "E42" '  -3.10%  '
Set-TextValue "E43" '  -1.96%  '
Set-TextValue "B44" 'Kaspa'
Set-TextValue "C44" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D44" '0.112'
Set-TextValue "E44" '  +2.03%  '
Set-TextValue "B45" 'Arweave'
Set-TextValue "C45" 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue "D45" '40.21'
Set-TextValue "E45" '  +2.30%  '
Set-TextValue "B46" 'VeChain'
Set-TextValue "C46" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D46" '0.0362'
Set-TextValue "E46" '  -3.12%  '
Set-TextValue "D47" '2.822.35'
Set-TextValue "E47" '  -1.64%  '
Set-TextValue "D48" '132.41'
Set-TextValue "E48" '  +1.01%  '
Set-TextValue "B49" 'USDe'
Set-TextValue "C49" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D49" '0.999'
Set-TextValue "E49" '  +0.01%  '
Set-TextValue "B50" 'InjectiveProtocol'
Set-TextValue "C50" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D50" '25.15'
Set-TextValue "E50" '  -1.01%  '
Set-TextValue "D51" '2.25'
Set-TextValue "E51" '  -2.14%  '
